$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the N1 header (missing closing parenthesis)
$ws.Range("N1").Value = "Land Price (R per Ha)"

# Add the new "Total Land Price (R)" column header in O1, matching the
# existing header formatting (bold, centered, bordered) by copying N1's format
$ws.Range("O1").Value = "Total Land Price (R)"
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

# Populate the Total Land Price (R) values for each data row
# (Total Land Price = Expected Total Land Required (L) * Land Price per Ha (N))
$ws.Range("O2").Value = 72247782.81229447
$ws.Range("O3").Value = 20782635.09843159
$ws.Range("O4").Value = 22593800.7240546
$ws.Range("O5").Value = 232581616.7967931
$ws.Range("O6").Value = 464676097.8206922
$ws.Range("O7").Value = 20853736.30525625
$ws.Range("O8").Value = 62319368.09302074
$ws.Range("O9").Value = 37765179.12383224
$ws.Range("O10").Value = 242233887.2037748
$ws.Range("O11").Value = 31183721.91838449
$ws.Range("O12").Value = 222864603.1357782
$ws.Range("O13").Value = 3779931.300469709
$ws.Range("O14").Value = 33320255.08831654
$ws.Range("O15").Value = 56768718.01137516
$ws.Range("O16").Value = 24373193.11248379
$ws.Range("O17").Value = 47322808.09180501
$ws.Range("O18").Value = 63706712.53960895
$ws.Range("O19").Value = 430793046.6636483
$ws.Range("O20").Value = 29484760.20079239
$ws.Range("O21").Value = 13662971.44586958
$ws.Range("O22").Value = 52989783.04911486
$ws.Range("O23").Value = 15037271.59408078
$ws.Range("O24").Value = 265968968.6051389
$ws.Range("O25").Value = 91725669.25305666
$ws.Range("O26").Value = 534393842.7625261
$ws.Range("O27").Value = 70222983.77083158
$ws.Range("O28").Value = 33296145.17918894
